$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: cohort_year 2024, period_index 1 -> num_customers 135 -> 136, retention_rate recalculated
$ws.Range("C36").Value = 136
$ws.Range("E36").Value = 136 / 1930

# Row 37: cohort_year 2025, period_index 0 -> num_customers 865 -> 870, cohort_size 865 -> 870
$ws.Range("C37").Value = 870
$ws.Range("D37").Value = 870
